$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4780
$ws.Range("J51").Value = 5453.222
$ws.Range("L51").Value = 5453.222
$ws.Range("N51").Value = -6421.222

$ws.Range("H76").Value = 3372.0908
$ws.Range("I76").Value = 3136.625
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3136.625
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -2821.625
$ws.Range("N76").Value = -4630

$ws.Range("H79").Value = 3372.0908
$ws.Range("I79").Value = 3136.625
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3136.625
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2044.625
$ws.Range("N79").Value = -6184

$ws.Range("H100").Value = 28572700
$ws.Range("I100").Value = 28572700
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 28572700
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -28572159
$ws.Range("N100").ClearContents()

$ws.Range("H113").Value = 19001.2
$ws.Range("J113").Value = 20001.5
$ws.Range("L113").Value = 20001.5
$ws.Range("N113").Value = -26509.5

$ws.Range("H124").Value = 42065.715
$ws.Range("J124").Value = 42065.715
$ws.Range("L124").Value = 42065.715
$ws.Range("N124").Value = -51885.715

$ws.Range("H132").Value = 40163252
$ws.Range("I132").Value = 45638480
$ws.Range("J132").Value = 11600
$ws.Range("K132").Value = 136915440
$ws.Range("L132").Value = 34800
$ws.Range("M132").Value = -136912910
$ws.Range("N132").Value = -39860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 701.25
$ws.Range("I2").Value = 658.5714
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 658.5714
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -545.5714
$ws.Range("N2").Value = -1226

$ws.Range("H32").Value = 5934.5
$ws.Range("I32").Value = 3929.3914
$ws.Range("K32").Value = 3929.3914
$ws.Range("M32").Value = -3642.3914

$ws.Range("H109").Value = 26050
$ws.Range("J109").Value = 26050
$ws.Range("L109").Value = 26050
$ws.Range("N109").Value = -28824

$ws.Range("H112").Value = 24890.436
$ws.Range("J112").Value = 24890.436
$ws.Range("L112").Value = 24890.436
$ws.Range("N112").Value = -27844.436

$ws.Range("H116").Value = 701.25
$ws.Range("I116").Value = 658.5714
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 658.5714
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1635.4286
$ws.Range("N116").Value = -5588

$ws.Range("H132").Value = 2786.805
$ws.Range("I132").Value = 2096.1936
$ws.Range("J132").Value = 4927.7
$ws.Range("K132").Value = 6288.5808
$ws.Range("L132").Value = 14783.1
$ws.Range("M132").Value = -3758.5808
$ws.Range("N132").Value = -19843.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 701.25
$ws.Range("I3").Value = 658.5714
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 658.5714
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -544.5714
$ws.Range("N3").Value = -1228

$ws.Range("H61").Value = 29107
$ws.Range("J61").Value = 29107
$ws.Range("L61").Value = 29107
$ws.Range("N61").Value = -29733

$ws.Range("H134").Value = 2741.0535
$ws.Range("I134").Value = 1633.3871
$ws.Range("J134").Value = 4114.56
$ws.Range("K134").Value = 4900.1613
$ws.Range("L134").Value = 12343.68
$ws.Range("M134").Value = -2365.1613
$ws.Range("N134").Value = -17413.68

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10102812
$ws.Range("I16").Value = 13889990
$ws.Range("K16").Value = 13889990
$ws.Range("M16").Value = -13889703

$ws.Range("H31").Value = 14709237
$ws.Range("I31").Value = 1898.625
$ws.Range("J31").Value = 50006850
$ws.Range("K31").Value = 1898.625
$ws.Range("L31").Value = 50006850
$ws.Range("M31").Value = -1603.625
$ws.Range("N31").Value = -50007440

$ws.Range("H34").Value = 14709237
$ws.Range("I34").Value = 1898.625
$ws.Range("J34").Value = 50006850
$ws.Range("K34").Value = 1898.625
$ws.Range("L34").Value = 50006850
$ws.Range("M34").Value = -1696.625
$ws.Range("N34").Value = -50007254

$ws.Range("H50").Value = 31514.6
$ws.Range("J50").Value = 31514.6
$ws.Range("L50").Value = 31514.6
$ws.Range("N50").Value = -32764.6

$ws.Range("H99").Value = 11768448
$ws.Range("I99").Value = 18183784
$ws.Range("K99").Value = 18183784
$ws.Range("M99").Value = -18182286

$ws.Range("H113").Value = 10102812
$ws.Range("I113").Value = 13889990
$ws.Range("K113").Value = 13889990
$ws.Range("M113").Value = -13887820

$ws.Range("H126").Value = 11768448
$ws.Range("I126").Value = 18183784
$ws.Range("K126").Value = 54551352
$ws.Range("M126").Value = -54548882

$ws.Range("H132").Value = 3137.838
$ws.Range("I132").Value = 1788.5555
$ws.Range("K132").Value = 5365.666499999999
$ws.Range("M132").Value = -2835.666499999999

$ws.Range("H134").Value = 5883.3213
$ws.Range("I134").Value = 7014.294
$ws.Range("J134").Value = 4135.4546
$ws.Range("K134").Value = 21042.882
$ws.Range("L134").Value = 12406.3638
$ws.Range("M134").Value = -18507.882
$ws.Range("N134").Value = -17476.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 6944.5454
$ws.Range("I112").Value = 7800
$ws.Range("J112").Value = 6903.8096
$ws.Range("K112").Value = 23400
$ws.Range("L112").Value = 20711.4288
$ws.Range("M112").Value = -22292
$ws.Range("N112").Value = -22927.4288

$ws.Range("H123").Value = 2547.1428
$ws.Range("I123").Value = 3107.5
$ws.Range("K123").Value = 9322.5
$ws.Range("M123").Value = -6872.5

$ws.Range("H129").Value = 2175.6428
$ws.Range("I129").Value = 2741.111
$ws.Range("J129").Value = 1157.8
$ws.Range("K129").Value = 8223.332999999999
$ws.Range("L129").Value = 3473.4
$ws.Range("M129").Value = -3223.332999999999
$ws.Range("N129").Value = -13473.4

$ws.Range("H131").Value = 7693323.5
$ws.Range("I131").Value = 100002270
$ws.Range("J131").Value = 911.05
$ws.Range("K131").Value = 300006810
$ws.Range("L131").Value = 2733.15
$ws.Range("M131").Value = -300001770
$ws.Range("N131").Value = -12813.15

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3003
$ws.Range("J7").Value = 3003
$ws.Range("L7").Value = 3003
$ws.Range("N7").Value = -3227

$ws.Range("H8").Value = 3003
$ws.Range("J8").Value = 3003
$ws.Range("L8").Value = 3003
$ws.Range("N8").Value = -3281

$ws.Range("H123").Value = 10910.444
$ws.Range("J123").Value = 10910.444
$ws.Range("L123").Value = 10910.444
$ws.Range("N123").Value = -15810.444

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1962.8636
$ws.Range("I46").Value = 1664.3572
$ws.Range("J46").Value = 2485.25
$ws.Range("K46").Value = 1664.3572
$ws.Range("L46").Value = 2485.25
$ws.Range("M46").Value = -1476.3572
$ws.Range("N46").Value = -2861.25

$ws.Range("H55").Value = 351.9375
$ws.Range("I55").Value = 249.66667
$ws.Range("J55").Value = 483.42856
$ws.Range("K55").Value = 249.66667
$ws.Range("L55").Value = 483.42856
$ws.Range("M55").Value = -76.66667000000001
$ws.Range("N55").Value = -829.4285600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10000000
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999886

$ws.Range("H11").Value = 7500000
$ws.Range("I11").Value = 7500000
$ws.Range("K11").Value = 7500000
$ws.Range("M11").Value = -7499858

$ws.Range("H107").Value = 925.06665
$ws.Range("I107").Value = 635.4545000000001
$ws.Range("J107").Value = 1721.5
$ws.Range("K107").Value = 1906.3635
$ws.Range("L107").Value = 5164.5
$ws.Range("M107").Value = 13.63649999999984
$ws.Range("N107").Value = -9004.5
